$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rebuild the "Components" checklist table.
# Column layout: A = section label; (B,C) / (D,E) / (F,G) / (H,I) are four
# two-column groups, each with a merged header in row 1 ("label" + "count"
# columns below it).
# ---------------------------------------------------------------------------

$ws.Range("A1").Value = "Components"

# --- Turn Table (B:C) ------------------------------------------------------
$ws.Range("B1").Value = "Turn Table"
$ws.Range("B2").Value = "Shaft";      $ws.Range("C2").Value = 0
$ws.Range("B3").Value = "Bearing";    $ws.Range("C3").Value = 0
$ws.Range("B4").Value = "Housing";    $ws.Range("C4").Value = 0
$ws.Range("B5").Value = "Motor ";     $ws.Range("C5").Value = 0
$ws.Range("B6").Value = "Connetion "; $ws.Range("C6").Value = 0
$ws.Range("B7").Value = "Table ";     $ws.Range("C7").Value = 0

# --- Elevator (D:E) ---------------------------------------------------------
$ws.Range("D1").Value = "Elevator "
$ws.Range("D2").Value = "Movement "

# --- End effector (F:G) -----------------------------------------------------
$ws.Range("F1").Value = "End effector  "
$ws.Range("F2").Value = "Servo Motor ";  $ws.Range("G2").Value = 1
$ws.Range("F3").Value = "Gear";          $ws.Range("G3").Value = 1
$ws.Range("F4").Value = "blades";        $ws.Range("G4").Value = 0
$ws.Range("F5").Value = "Housing";       $ws.Range("G5").Value = 0
$ws.Range("F6").Value = "Bearing";       $ws.Range("G6").Value = 0
$ws.Range("F7").Value = "Tube";          $ws.Range("G7").Value = 0
$ws.Range("F8").Value = "Connection ";   $ws.Range("G8").Value = 0
$ws.Range("F9").Value = "Shaft";         $ws.Range("G9").Value = 0
$ws.Range("F10").Value = "Scissor ";     $ws.Range("G10").Value = 0

# --- Arm (H:I) ---------------------------------------------------------------
$ws.Range("H1").Value = "Arm "
$ws.Range("H2").Value = "Beam ";        $ws.Range("I2").Value = 0
$ws.Range("H3").Value = "Power Screw";  $ws.Range("I3").Value = 0
$ws.Range("H4").Value = "Motor";        $ws.Range("I4").Value = 0
$ws.Range("H5").Value = "Housing";      $ws.Range("I5").Value = 0
$ws.Range("H6").Value = "Bearing";      $ws.Range("I6").Value = 0

# ---------------------------------------------------------------------------
# Header formatting: center-align each header pair, then merge it.
# ---------------------------------------------------------------------------
$ws.Range("B1:C1").HorizontalAlignment = -4108
$ws.Range("D1:E1").HorizontalAlignment = -4108
$ws.Range("F1:G1").HorizontalAlignment = -4108
$ws.Range("H1:I1").HorizontalAlignment = -4108

$ws.Range("B1:C1").Merge()
$ws.Range("D1:E1").Merge()
$ws.Range("F1:G1").Merge()
$ws.Range("H1:I1").Merge()

# ---------------------------------------------------------------------------
# Column widths / zoom / selection
# ---------------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 11.8
$ws.Columns("H").ColumnWidth = 13.8

$excel.ActiveWindow.Zoom = 213
$ws.Range("G18").Select()
